$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 243; existing rows 243-291 shift down to 244-292.
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with the new weekly record.
$ws.Cells.Item(243, 1).Value = 5
$ws.Cells.Item(243, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(243, 3).Value = "Maule"
$ws.Cells.Item(243, 4).Value = 44637
$ws.Cells.Item(243, 5).Value = 7
$ws.Cells.Item(243, 6).Value = 100112023
$ws.Cells.Item(243, 7).Value = "Brócoli"
$ws.Cells.Item(243, 8).Value = "Sin especificar"
$ws.Cells.Item(243, 9).Value = "Primera"
$ws.Cells.Item(243, 10).Value = 5000
$ws.Cells.Item(243, 11).Value = 400
$ws.Cells.Item(243, 12).Value = 400
$ws.Cells.Item(243, 13).Value = 400
$ws.Cells.Item(243, 14).Value = "$/unidad"
$ws.Cells.Item(243, 15).Value = "Región del Maule"
$ws.Cells.Item(243, 16).Value = 400
$ws.Cells.Item(243, 17).Value = 1
$ws.Cells.Item(243, 18).Value = "Hortaliza"
